$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "[2023-01-01]`n内膜組織診の結果、漿液性卵巣癌、Stage IIIAと診断。 抗凝固薬（エリキュース）服用中`n`n[2023-01-26]`n治療方針：腹腔鏡下子宮全摘+両側付属器切除術の方針。 腎機能低下（Cr 1.9 mg/dl）`n`n[2023-02-20]`n手術記録：腹腔鏡下子宮全摘施行。手術時間242分、出血量275ml。`n`n[2023-03-15]`n術後化学療法としてパクリタキセル+カルボプラチンを開始。`n`n[2023-03-24]`n術後化学療法としてweekly PTXを開始。`n`n[2023-04-15]`n術後化学療法としてPTX+CBDCAを開始。"
$ws.Range("C2").Value = "漿液性卵巣癌"
$ws.Range("D2").Value = "カルテ1行目に記載された最新の診断名"
$ws.Range("E2").Value = "Stage IIIA"
$ws.Range("F2").Value = "2023-01-01のカルテに記載された最新のステージ情報"
$ws.Range("G2").Value = "2023/01/01 内膜組織診: 漿液性卵巣癌、Stage IIIA"
$ws.Range("H2").Value = "2023-01-01の記載から内膜組織診の結果が抽出されました。"
$ws.Range("B3").Value = "[2023-01-01]`n組織診断の結果、卵巣腫瘍、Stage IICと診断。 高度肥満（BMI 36.2）`n`n[2023-01-22]`n治療方針：開腹子宮全摘出術の方針。`n`n[2023-02-11]`n手術記録：腹腔鏡下子宮全摘施行。手術時間255分、出血量853ml。`n`n[2023-03-04]`n術後化学療法としてドセタキセル+カルボプラチンを開始。`n`n[2023-04-02]`n術後化学療法としてドセタキセル+カルボプラチンを開始。"
$ws.Range("C3").Value = "{`"result`": `"卵巣腫瘍`", `"reason`": `"2023-01-01のカルテに記載された最新の診断名`"}`nただし、`"卵巣腫瘍`"はがんの具体的な診断名とは言えません。病理診断などで具体的な組織型が記載されていないため、より具体的な診断名が見つからなかった場合、以下のように出力します。`n{`"result`": `"記載なし`", `"reason`": `"がんの具体的な診断名の記載が見つかりませんでした`"}"
$ws.Range("D3").Value = "JSONエラー"
$ws.Range("E3").Value = "Stage IIC"
$ws.Range("F3").Value = "2023-01-01の組織診断結果に記載された最新のステージ情報"
$ws.Range("G3").Value = "2023/02/11 腹腔鏡下子宮全摘施行"
$ws.Range("H3").Value = "2023-02-11の手術記録から抽出"
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "[2023-01-01]`n組織診断の結果、卵巣癌、Stage IBと診断。`n`n[2023-01-09]`n治療方針：腹腔鏡下全摘術の方針。 抗凝固薬（イグザレルト）服用中`n`n[2023-02-03]`n手術記録：開腹子宮全摘出術施行。手術時間226分、出血量72ml。`n`n[2023-03-03]`n術後化学療法としてドセタキセル+カルボプラチンを開始。`n`n[2023-03-30]`n術後化学療法としてTC療法を開始。`n`n[2023-04-29]`n術後化学療法としてweekly パクリタキセルを開始。 心機能低下（EF 47%）あり"
$ws.Range("C4").Value = "卵巣癌"
$ws.Range("D4").Value = "2023-01-01のカルテに記載された最新の診断名"
$ws.Range("E4").Value = "Stage IB"
$ws.Range("F4").Value = "2023-01-01の組織診断の結果に記載された最新のステージ情報"
$ws.Range("G4").Value = "2023/02/03 開腹子宮全摘出術"
$ws.Range("H4").Value = "2023-02-03の手術記録から抽出"

# Re-fit row heights so no explicit custom row height is persisted
# (multi-line text assignment otherwise bakes in an auto row height).
$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
$ws.Rows(4).AutoFit()
